# Update Leve profit metrics (currentAveragePrice*, LevePrice*, LeveProfit*) columns
# with freshly pulled market-board figures, per sheet/row, as produced by the
# scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 427.1
$ws.Cells.Item(32, 9).Value = 324.14285
$ws.Cells.Item(32, 10).Value = 667.3333
$ws.Cells.Item(32, 11).Value = 324.14285
$ws.Cells.Item(32, 12).Value = 667.3333
$ws.Cells.Item(32, 13).Value = 1.85714999999999
$ws.Cells.Item(32, 14).Value = -1319.3333
$ws.Cells.Item(74, 8).Value = 4723.625
$ws.Cells.Item(74, 9).Value = 3881.6667
$ws.Cells.Item(74, 11).Value = 3881.6667
$ws.Cells.Item(74, 13).Value = -2945.6667
$ws.Cells.Item(77, 8).Value = 4723.625
$ws.Cells.Item(77, 9).Value = 3881.6667
$ws.Cells.Item(77, 11).Value = 19408.3335
$ws.Cells.Item(77, 13).Value = -14728.3335
$ws.Cells.Item(132, 8).Value = 2015.1724
$ws.Cells.Item(132, 9).Value = 1682.7407
$ws.Cells.Item(132, 10).Value = 6503
$ws.Cells.Item(132, 11).Value = 5048.2221
$ws.Cells.Item(132, 12).Value = 19509
$ws.Cells.Item(132, 13).Value = -2518.2221
$ws.Cells.Item(132, 14).Value = -24569
$ws.Cells.Item(135, 8).Value = 57693548
$ws.Cells.Item(135, 9).Value = 20834638
$ws.Cells.Item(135, 10).Value = 500000500
$ws.Cells.Item(135, 11).Value = 187511742
$ws.Cells.Item(135, 12).Value = 4500004500
$ws.Cells.Item(135, 13).Value = -187509207
$ws.Cells.Item(135, 14).Value = -4500009570

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2145.3125
$ws.Cells.Item(2, 9).Value = 2208.3333
$ws.Cells.Item(2, 10).Value = 1200
$ws.Cells.Item(2, 11).Value = 2208.3333
$ws.Cells.Item(2, 12).Value = 1200
$ws.Cells.Item(2, 13).Value = -2095.3333
$ws.Cells.Item(2, 14).Value = -1426
$ws.Cells.Item(32, 8).Value = 15958.014
$ws.Cells.Item(32, 9).Value = 17736.613
$ws.Cells.Item(32, 10).Value = 4930.7
$ws.Cells.Item(32, 11).Value = 17736.613
$ws.Cells.Item(32, 12).Value = 4930.7
$ws.Cells.Item(32, 13).Value = -17449.613
$ws.Cells.Item(32, 14).Value = -5504.7
$ws.Cells.Item(61, 8).Value = 7748.3403
$ws.Cells.Item(61, 9).Value = 4266.081
$ws.Cells.Item(61, 11).Value = 4266.081
$ws.Cells.Item(61, 13).Value = -4054.081
$ws.Cells.Item(74, 8).Value = 1989.75
$ws.Cells.Item(74, 9).Value = 2013.0344
$ws.Cells.Item(74, 10).Value = 1764.6666
$ws.Cells.Item(74, 11).Value = 2013.0344
$ws.Cells.Item(74, 12).Value = 1764.6666
$ws.Cells.Item(74, 13).Value = -1139.0344
$ws.Cells.Item(74, 14).Value = -3512.6666
$ws.Cells.Item(77, 8).Value = 1989.75
$ws.Cells.Item(77, 9).Value = 2013.0344
$ws.Cells.Item(77, 10).Value = 1764.6666
$ws.Cells.Item(77, 11).Value = 10065.172
$ws.Cells.Item(77, 12).Value = 8823.333000000001
$ws.Cells.Item(77, 13).Value = -5697.172
$ws.Cells.Item(77, 14).Value = -17559.333
$ws.Cells.Item(116, 8).Value = 2145.3125
$ws.Cells.Item(116, 9).Value = 2208.3333
$ws.Cells.Item(116, 10).Value = 1200
$ws.Cells.Item(116, 11).Value = 2208.3333
$ws.Cells.Item(116, 12).Value = 1200
$ws.Cells.Item(116, 13).Value = 85.66670000000022
$ws.Cells.Item(116, 14).Value = -5788
$ws.Cells.Item(132, 8).Value = 1645.4762
$ws.Cells.Item(132, 9).Value = 1331.875
$ws.Cells.Item(132, 10).Value = 2649
$ws.Cells.Item(132, 11).Value = 3995.625
$ws.Cells.Item(132, 12).Value = 7947
$ws.Cells.Item(132, 13).Value = -1465.625
$ws.Cells.Item(132, 14).Value = -13007
$ws.Cells.Item(136, 8).Value = 7748.3403
$ws.Cells.Item(136, 9).Value = 4266.081
$ws.Cells.Item(136, 11).Value = 12798.243
$ws.Cells.Item(136, 13).Value = -10248.243
$ws.Cells.Item(138, 8).Value = 63198.332
$ws.Cells.Item(138, 10).Value = 63198.332
$ws.Cells.Item(138, 12).Value = 63198.332
$ws.Cells.Item(138, 14).Value = -73478.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2145.3125
$ws.Cells.Item(3, 9).Value = 2208.3333
$ws.Cells.Item(3, 10).Value = 1200
$ws.Cells.Item(3, 11).Value = 2208.3333
$ws.Cells.Item(3, 12).Value = 1200
$ws.Cells.Item(3, 13).Value = -2094.3333
$ws.Cells.Item(3, 14).Value = -1428
$ws.Cells.Item(80, 8).Value = 219
$ws.Cells.Item(80, 10).Value = 243.71428
$ws.Cells.Item(80, 12).Value = 243.71428
$ws.Cells.Item(80, 14).Value = -2239.71428
$ws.Cells.Item(83, 8).Value = 219
$ws.Cells.Item(83, 10).Value = 243.71428
$ws.Cells.Item(83, 12).Value = 1218.5714
$ws.Cells.Item(83, 14).Value = -11202.5714
$ws.Cells.Item(94, 8).Value = 1457.8
$ws.Cells.Item(94, 9).Value = 1238.0834
$ws.Cells.Item(94, 10).Value = 2336.6667
$ws.Cells.Item(94, 11).Value = 1238.0834
$ws.Cells.Item(94, 12).Value = 2336.6667
$ws.Cells.Item(94, 13).Value = -787.0834
$ws.Cells.Item(94, 14).Value = -3238.6667
$ws.Cells.Item(107, 8).Value = 1542.2
$ws.Cells.Item(107, 9).Value = 905.5
$ws.Cells.Item(107, 10).Value = 1966.6666
$ws.Cells.Item(107, 11).Value = 905.5
$ws.Cells.Item(107, 12).Value = 1966.6666
$ws.Cells.Item(107, 13).Value = 1014.5
$ws.Cells.Item(107, 14).Value = -5806.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1289.4286
$ws.Cells.Item(99, 9).Value = 1162.5714
$ws.Cells.Item(99, 11).Value = 1162.5714
$ws.Cells.Item(99, 13).Value = 335.4286
$ws.Cells.Item(126, 8).Value = 1289.4286
$ws.Cells.Item(126, 9).Value = 1162.5714
$ws.Cells.Item(126, 11).Value = 3487.7142
$ws.Cells.Item(126, 13).Value = -1017.7142
$ws.Cells.Item(138, 8).Value = 62338.75
$ws.Cells.Item(138, 10).Value = 62338.75
$ws.Cells.Item(138, 12).Value = 62338.75
$ws.Cells.Item(138, 14).Value = -72618.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 499
$ws.Cells.Item(18, 9).Value = 521.1111
$ws.Cells.Item(18, 10).Value = 300
$ws.Cells.Item(18, 11).Value = 1563.3333
$ws.Cells.Item(18, 12).Value = 900
$ws.Cells.Item(18, 13).Value = -1394.3333
$ws.Cells.Item(18, 14).Value = -1238
$ws.Cells.Item(97, 8).Value = 3000.5
$ws.Cells.Item(97, 9).Value = 501
$ws.Cells.Item(97, 10).Value = 5500
$ws.Cells.Item(97, 11).Value = 1503
$ws.Cells.Item(97, 12).Value = 16500
$ws.Cells.Item(97, 13).Value = -1007
$ws.Cells.Item(97, 14).Value = -17492

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 4101.2964
$ws.Cells.Item(102, 9).Value = 3590.6316
$ws.Cells.Item(102, 10).Value = 5314.125
$ws.Cells.Item(102, 11).Value = 3590.6316
$ws.Cells.Item(102, 12).Value = 5314.125
$ws.Cells.Item(102, 13).Value = -1968.6316
$ws.Cells.Item(102, 14).Value = -8558.125
$ws.Cells.Item(132, 8).Value = 7689.476
$ws.Cells.Item(132, 9).Value = 2807.4443
$ws.Cells.Item(132, 10).Value = 11351
$ws.Cells.Item(132, 11).Value = 8422.332900000001
$ws.Cells.Item(132, 12).Value = 34053
$ws.Cells.Item(132, 13).Value = -5892.332900000001
$ws.Cells.Item(132, 14).Value = -39113

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 718.75
$ws.Cells.Item(16, 9).Value = 575
$ws.Cells.Item(16, 11).Value = 575
$ws.Cells.Item(16, 13).Value = -405
$ws.Cells.Item(82, 8).Value = 1449.4
$ws.Cells.Item(82, 9).Value = 1111.75
$ws.Cells.Item(82, 10).Value = 2800
$ws.Cells.Item(82, 11).Value = 1111.75
$ws.Cells.Item(82, 12).Value = 2800
$ws.Cells.Item(82, 13).Value = -750.75
$ws.Cells.Item(82, 14).Value = -3522
$ws.Cells.Item(85, 8).Value = 1449.4
$ws.Cells.Item(85, 9).Value = 1111.75
$ws.Cells.Item(85, 10).Value = 2800
$ws.Cells.Item(85, 11).Value = 1111.75
$ws.Cells.Item(85, 12).Value = 2800
$ws.Cells.Item(85, 13).Value = 136.25
$ws.Cells.Item(85, 14).Value = -5296
$ws.Cells.Item(87, 8).Value = 30839.5
$ws.Cells.Item(87, 9).Value = 7980
$ws.Cells.Item(87, 11).Value = 7980
$ws.Cells.Item(87, 13).Value = -6857
$ws.Cells.Item(88, 8).Value = 40162
$ws.Cells.Item(88, 10).Value = 40189
$ws.Cells.Item(88, 12).Value = 40189
$ws.Cells.Item(88, 14).Value = -41045
$ws.Cells.Item(90, 8).Value = 30839.5
$ws.Cells.Item(90, 9).Value = 7980
$ws.Cells.Item(90, 11).Value = 23940
$ws.Cells.Item(90, 13).Value = -18324
$ws.Cells.Item(91, 8).Value = 40162
$ws.Cells.Item(91, 10).Value = 40189
$ws.Cells.Item(91, 12).Value = 40189
$ws.Cells.Item(91, 14).Value = -43153
$ws.Cells.Item(132, 8).Value = 4089
$ws.Cells.Item(132, 9).Value = 3661.0625
$ws.Cells.Item(132, 10).Value = 4849.778
$ws.Cells.Item(132, 11).Value = 10983.1875
$ws.Cells.Item(132, 12).Value = 14549.334
$ws.Cells.Item(132, 13).Value = -8453.1875
$ws.Cells.Item(132, 14).Value = -19609.334
$ws.Cells.Item(138, 8).Value = 40000
$ws.Cells.Item(138, 10).Value = 40000
$ws.Cells.Item(138, 12).Value = 40000
$ws.Cells.Item(138, 14).Value = -50280

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 15909.363
$ws.Cells.Item(2, 9).Value = 28167.334
$ws.Cells.Item(2, 10).Value = 11312.625
$ws.Cells.Item(2, 11).Value = 28167.334
$ws.Cells.Item(2, 12).Value = 11312.625
$ws.Cells.Item(2, 13).Value = -28055.334
$ws.Cells.Item(2, 14).Value = -11536.625
$ws.Cells.Item(56, 8).Value = 27970.334
$ws.Cells.Item(56, 9).Value = 7285
$ws.Cells.Item(56, 10).Value = 38313
$ws.Cells.Item(56, 11).Value = 7285
$ws.Cells.Item(56, 12).Value = 38313
$ws.Cells.Item(56, 13).Value = -6571
$ws.Cells.Item(56, 14).Value = -39741
$ws.Cells.Item(107, 8).Value = 3296.8333
$ws.Cells.Item(107, 10).Value = 6098.6665
$ws.Cells.Item(107, 12).Value = 18295.9995
$ws.Cells.Item(107, 14).Value = -22135.9995
